$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, shifting existing rows 13-28 down to 14-29.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new weekly data point.
$ws.Cells.Item(13, 1).Value = 10
$ws.Cells.Item(13, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(13, 3).Value = "La Araucanía"
$ws.Cells.Item(13, 4).Value = 44662
$ws.Cells.Item(13, 5).Value = 9
$ws.Cells.Item(13, 6).Value = "Fruta"
$ws.Cells.Item(13, 7).Value = 100108
$ws.Cells.Item(13, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(13, 9).Value = 100108003
$ws.Cells.Item(13, 10).Value = "Maracuyá"
$ws.Cells.Item(13, 11).Value = "Sin especificar"
$ws.Cells.Item(13, 12).Value = "Primera"
$ws.Cells.Item(13, 13).Value = 15
$ws.Cells.Item(13, 14).Value = 30000
$ws.Cells.Item(13, 15).Value = 30000
$ws.Cells.Item(13, 16).Value = 30000
$ws.Cells.Item(13, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(13, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(13, 19).Value = 1667
$ws.Cells.Item(13, 20).Value = 18
